$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32. This shifts the existing rows 32:67
# (and all their formatting, e.g. the date number format on column D)
# down to rows 33:68, matching the weekly update reflected in the diff.
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with this week's new price record.
$ws.Range("A32").Value = 8
$ws.Range("B32").Value = "Terminal La Palmera de La Serena"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = 44601
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = 100112030
$ws.Range("G32").Value = "Poroto granado"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 300
$ws.Range("K32").Value = 31000
$ws.Range("L32").Value = 32000
$ws.Range("M32").Value = 31500
$ws.Range("N32").Value = "$/malla 25 kilos"
$ws.Range("O32").Value = "Provincia del Elquí"
$ws.Range("P32").Value = 1260
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"
